$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Cells.Item(15, 8).Value = 2098.3389
$ws_ALC.Cells.Item(15, 9).Value = 2098.3389
$ws_ALC.Cells.Item(15, 11).Value = 6295.0167
$ws_ALC.Cells.Item(15, 13).Value = -6126.0167
$ws_ALC.Cells.Item(112, 8).Value = 1583.2
$ws_ALC.Cells.Item(112, 9).Value = 3262.5
$ws_ALC.Cells.Item(112, 10).Value = 1163.375
$ws_ALC.Cells.Item(112, 11).Value = 9787.5
$ws_ALC.Cells.Item(112, 12).Value = 3490.125
$ws_ALC.Cells.Item(112, 13).Value = -8679.5
$ws_ALC.Cells.Item(112, 14).Value = -5706.125
$ws_ALC.Cells.Item(113, 8).Value = 5653.3335
$ws_ALC.Cells.Item(113, 9).Value = 4784
$ws_ALC.Cells.Item(113, 10).Value = 10000
$ws_ALC.Cells.Item(113, 11).Value = 4784
$ws_ALC.Cells.Item(113, 12).Value = 10000
$ws_ALC.Cells.Item(113, 13).Value = -1530
$ws_ALC.Cells.Item(113, 14).Value = -16508
$ws_ALC.Cells.Item(129, 8).Value = 904.57733
$ws_ALC.Cells.Item(129, 9).Value = 1459.4
$ws_ALC.Cells.Item(129, 10).Value = 874.4239
$ws_ALC.Cells.Item(129, 11).Value = 4378.200000000001
$ws_ALC.Cells.Item(129, 12).Value = 2623.2717
$ws_ALC.Cells.Item(129, 13).Value = 621.7999999999993
$ws_ALC.Cells.Item(129, 14).Value = -12623.2717
$ws_ALC.Cells.Item(132, 8).Value = 6905873
$ws_ALC.Cells.Item(132, 9).Value = 8343681.5
$ws_ALC.Cells.Item(132, 10).Value = 4393.2
$ws_ALC.Cells.Item(132, 11).Value = 25031044.5
$ws_ALC.Cells.Item(132, 12).Value = 13179.6
$ws_ALC.Cells.Item(132, 13).Value = -25028514.5
$ws_ALC.Cells.Item(132, 14).Value = -18239.6
$ws_ALC.Cells.Item(133, 8).Value = 20719.215
$ws_ALC.Cells.Item(133, 10).Value = 20719.215
$ws_ALC.Cells.Item(133, 12).Value = 20719.215
$ws_ALC.Cells.Item(133, 14).Value = -30839.215

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Cells.Item(45, 8).Value = 1472.7878
$ws_ARM.Cells.Item(45, 9).Value = 1096.48
$ws_ARM.Cells.Item(45, 10).Value = 2648.75
$ws_ARM.Cells.Item(45, 11).Value = 1096.48
$ws_ARM.Cells.Item(45, 12).Value = 2648.75
$ws_ARM.Cells.Item(45, 13).Value = -719.48
$ws_ARM.Cells.Item(45, 14).Value = -3402.75
$ws_ARM.Cells.Item(97, 8).Value = 2010
$ws_ARM.Cells.Item(97, 9).Value = 2010
$ws_ARM.Cells.Item(97, 10).Value = 0
$ws_ARM.Cells.Item(97, 11).Value = 2010
$ws_ARM.Cells.Item(97, 12).Value = 0
$ws_ARM.Cells.Item(97, 13).Value = -1514
$ws_ARM.Cells.Item(97, 14).ClearContents()
$ws_ARM.Cells.Item(132, 8).Value = 2430.6206
$ws_ARM.Cells.Item(132, 9).Value = 2212.6
$ws_ARM.Cells.Item(132, 10).Value = 3185.3076
$ws_ARM.Cells.Item(132, 11).Value = 6637.799999999999
$ws_ARM.Cells.Item(132, 12).Value = 9555.9228
$ws_ARM.Cells.Item(132, 13).Value = -4107.799999999999
$ws_ARM.Cells.Item(132, 14).Value = -14615.9228

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Cells.Item(99, 8).Value = 3739.5
$ws_BSM.Cells.Item(99, 9).Value = 2640.3333
$ws_BSM.Cells.Item(99, 10).Value = 4563.875
$ws_BSM.Cells.Item(99, 11).Value = 2640.3333
$ws_BSM.Cells.Item(99, 12).Value = 4563.875
$ws_BSM.Cells.Item(99, 13).Value = -1142.3333
$ws_BSM.Cells.Item(99, 14).Value = -7559.875
$ws_BSM.Cells.Item(105, 8).Value = 1722.963
$ws_BSM.Cells.Item(105, 9).Value = 1615.5
$ws_BSM.Cells.Item(105, 10).Value = 2030
$ws_BSM.Cells.Item(105, 11).Value = 1615.5
$ws_BSM.Cells.Item(105, 12).Value = 2030
$ws_BSM.Cells.Item(105, 13).Value = 131.5
$ws_BSM.Cells.Item(105, 14).Value = -5524

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Cells.Item(16, 8).Value = 2183.4375
$ws_CRP.Cells.Item(16, 9).Value = 1512.375
$ws_CRP.Cells.Item(16, 10).Value = 2854.5
$ws_CRP.Cells.Item(16, 11).Value = 1512.375
$ws_CRP.Cells.Item(16, 12).Value = 2854.5
$ws_CRP.Cells.Item(16, 13).Value = -1225.375
$ws_CRP.Cells.Item(16, 14).Value = -3428.5
$ws_CRP.Cells.Item(31, 8).Value = 4481.3335
$ws_CRP.Cells.Item(31, 9).Value = 3259.4666
$ws_CRP.Cells.Item(31, 10).Value = 6008.6665
$ws_CRP.Cells.Item(31, 11).Value = 3259.4666
$ws_CRP.Cells.Item(31, 12).Value = 6008.6665
$ws_CRP.Cells.Item(31, 13).Value = -2964.4666
$ws_CRP.Cells.Item(31, 14).Value = -6598.6665
$ws_CRP.Cells.Item(34, 8).Value = 4481.3335
$ws_CRP.Cells.Item(34, 9).Value = 3259.4666
$ws_CRP.Cells.Item(34, 10).Value = 6008.6665
$ws_CRP.Cells.Item(34, 11).Value = 3259.4666
$ws_CRP.Cells.Item(34, 12).Value = 6008.6665
$ws_CRP.Cells.Item(34, 13).Value = -3057.4666
$ws_CRP.Cells.Item(34, 14).Value = -6412.6665
$ws_CRP.Cells.Item(52, 8).Value = 67500
$ws_CRP.Cells.Item(52, 10).Value = 67500
$ws_CRP.Cells.Item(52, 12).Value = 67500
$ws_CRP.Cells.Item(52, 14).Value = -68088
$ws_CRP.Cells.Item(107, 8).Value = 1556.3
$ws_CRP.Cells.Item(107, 9).Value = 607.1429000000001
$ws_CRP.Cells.Item(107, 10).Value = 3771
$ws_CRP.Cells.Item(107, 11).Value = 607.1429000000001
$ws_CRP.Cells.Item(107, 12).Value = 3771
$ws_CRP.Cells.Item(107, 13).Value = 1312.8571
$ws_CRP.Cells.Item(107, 14).Value = -7611
$ws_CRP.Cells.Item(113, 8).Value = 2183.4375
$ws_CRP.Cells.Item(113, 9).Value = 1512.375
$ws_CRP.Cells.Item(113, 10).Value = 2854.5
$ws_CRP.Cells.Item(113, 11).Value = 1512.375
$ws_CRP.Cells.Item(113, 12).Value = 2854.5
$ws_CRP.Cells.Item(113, 13).Value = 657.625
$ws_CRP.Cells.Item(113, 14).Value = -7194.5

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Cells.Item(102, 8).Value = 31842.656
$ws_GSM.Cells.Item(102, 9).Value = 2277.1667
$ws_GSM.Cells.Item(102, 10).Value = 47268.13
$ws_GSM.Cells.Item(102, 11).Value = 2277.1667
$ws_GSM.Cells.Item(102, 12).Value = 47268.13
$ws_GSM.Cells.Item(102, 13).Value = -655.1667000000002
$ws_GSM.Cells.Item(102, 14).Value = -50512.13

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Cells.Item(14, 8).Value = 202159.23
$ws_LTW.Cells.Item(14, 9).Value = 1250500
$ws_LTW.Cells.Item(14, 11).Value = 1250500
$ws_LTW.Cells.Item(14, 13).Value = -1250328
$ws_LTW.Cells.Item(22, 8).Value = 250000370
$ws_LTW.Cells.Item(22, 9).Value = 500000260
$ws_LTW.Cells.Item(22, 11).Value = 500000260
$ws_LTW.Cells.Item(22, 13).Value = -499999965
$ws_LTW.Cells.Item(27, 8).Value = 250000370
$ws_LTW.Cells.Item(27, 9).Value = 500000260
$ws_LTW.Cells.Item(27, 11).Value = 500000260
$ws_LTW.Cells.Item(27, 13).Value = -500000153
$ws_LTW.Cells.Item(132, 8).Value = 2511.0881
$ws_LTW.Cells.Item(132, 9).Value = 1591.4231
$ws_LTW.Cells.Item(132, 10).Value = 5500
$ws_LTW.Cells.Item(132, 11).Value = 4774.2693
$ws_LTW.Cells.Item(132, 12).Value = 16500
$ws_LTW.Cells.Item(132, 13).Value = -2244.2693
$ws_LTW.Cells.Item(132, 14).Value = -21560
$ws_LTW.Cells.Item(137, 8).Value = 44900
$ws_LTW.Cells.Item(137, 10).Value = 44900
$ws_LTW.Cells.Item(137, 12).Value = 44900
$ws_LTW.Cells.Item(137, 14).Value = -55100
$ws_LTW.Cells.Item(139, 8).Value = 29750
$ws_LTW.Cells.Item(139, 10).Value = 29750
$ws_LTW.Cells.Item(139, 12).Value = 29750
$ws_LTW.Cells.Item(139, 14).Value = -40030
$ws_LTW.Cells.Item(141, 8).Value = 29583.334
$ws_LTW.Cells.Item(141, 10).Value = 29583.334
$ws_LTW.Cells.Item(141, 12).Value = 29583.334
$ws_LTW.Cells.Item(141, 14).Value = -39943.334

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Cells.Item(45, 8).Value = 0
$ws_WVR.Cells.Item(45, 10).Value = 0
$ws_WVR.Cells.Item(45, 12).Value = 0
$ws_WVR.Cells.Item(45, 14).ClearContents()
$ws_WVR.Cells.Item(132, 8).Value = 19393.451
$ws_WVR.Cells.Item(132, 9).Value = 3056.3809
$ws_WVR.Cells.Item(132, 10).Value = 53701.3
$ws_WVR.Cells.Item(132, 11).Value = 9169.1427
$ws_WVR.Cells.Item(132, 12).Value = 161103.9
$ws_WVR.Cells.Item(132, 13).Value = -6639.1427
$ws_WVR.Cells.Item(132, 14).Value = -166163.9
$ws_WVR.Cells.Item(135, 8).Value = 150715
$ws_WVR.Cells.Item(135, 10).Value = 150715
$ws_WVR.Cells.Item(135, 12).Value = 150715
$ws_WVR.Cells.Item(135, 14).Value = -160855
$ws_WVR.Cells.Item(136, 8).Value = 1347.2285
$ws_WVR.Cells.Item(136, 9).Value = 756.375
$ws_WVR.Cells.Item(136, 10).Value = 2636.3635
$ws_WVR.Cells.Item(136, 11).Value = 2269.125
$ws_WVR.Cells.Item(136, 12).Value = 7909.0905
$ws_WVR.Cells.Item(136, 13).Value = 280.875
$ws_WVR.Cells.Item(136, 14).Value = -13009.0905
$ws_WVR.Cells.Item(137, 8).Value = 66900
$ws_WVR.Cells.Item(137, 10).Value = 66900
$ws_WVR.Cells.Item(137, 12).Value = 66900
$ws_WVR.Cells.Item(137, 14).Value = -77100
$ws_WVR.Cells.Item(141, 8).Value = 29642.857
$ws_WVR.Cells.Item(141, 10).Value = 29642.857
$ws_WVR.Cells.Item(141, 12).Value = 29642.857
$ws_WVR.Cells.Item(141, 14).Value = -40002.857

